# update: 1/12/2025: Update last 30 days report
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (Silverfort / Nordics RSM / Jesper Damm-Skogh) moved from "CV Sent" to "1st Interview"
$ws.Range("E22").Value = "1st Interview"
$ws.Range("F22").Value = 45991

# Insert a new row 23 for the second Silverfort / Nordics RSM candidate, pushing
# the existing row 23 (Blockaid / SDR Manager / Tae Kim) down to row 24.
$ws.Rows.Item(23).Insert()

$ws.Range("A23").Value = 820
$ws.Range("B23").Value = "Silverfort"
$ws.Range("C23").Value = "Nordics RSM"
$ws.Range("D23").Value = "Marc Solis"
$ws.Range("E23").Value = "1st Interview"
$ws.Range("F23").Value = 45991
